# Apply the commit's edits:
#   1. Header-row heights in the 5 ranova tables go from 615/571 twips to 637 twips
#      (637 twips = 31.85 points).
#   2. The Greek "χ" (chi) character in each table's "χ2" header cell is replaced
#      with the mis-encoded byte sequence "Ï‡" (U+00CF, U+2021) - this mirrors the
#      exact mojibake text that shows up in the author's commit.

$d = $word.ActiveDocument

# --- 1. Fix up the header row height on every table in the document ---
# 637 twentieths-of-a-point (dxa) == 31.85 points; Word's Row.Height is in points.
$newHeightPoints = 637 / 20
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $headerRow = $tbl.Rows.Item(1)
    $headerRow.Height = $newHeightPoints
}

# --- 2. Replace every "χ" with the mojibake "Ï‡" (Ã¯ + double-dagger) ---
$chi = [string][char]0x3C7
$mojibake = [string]::Concat([char]0xCF, [char]0x2021)

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute($chi, $true, $false, $false, $false, $false, $true, 1, $false, $mojibake, 2)
